# "Script to download new sets finished": Wizards Play Network 2008 (PWPN)
# now has its full set name on A1 and its card list cleared out (A2 is
# blanked, the old A3 row is removed) so the next run can repopulate it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1: set title becomes the full set name.
$ws.Range("A1").Value = "Wizards Play Network 2008 (PWPN)"

# Drop the old card rows (A2: "Sprouting Thrinax", A3: "Woolly Thoctar").
# Clear A3 first so it is removed entirely, shrinking the sheet to A1:A2.
$ws.Range("A3").ClearContents()

# A2 becomes blank but stays part of the used range - clear its contents,
# then nudge its style so the row is retained as an empty cell rather than
# being pruned along with A3.
$ws.Range("A2").ClearContents()
$ws.Range("A2").Style = "Normal"
